$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D (average) and E-J
$ws.Range("D2").Value = 2010612.8556717001
$ws.Range("E2").Value = 0.53000000000000003
$ws.Range("F2").Value = 1408.8077184772701
$ws.Range("G2").Value = 0.99999988822061503
$ws.Range("H2").Value = 4277899.2147375103
$ws.Range("I2").Value = 23976.794391179901
$ws.Range("J2").Value = 3036.5387402628298

# Rows 3-63: update E-J only (D stays blank)
# Row 3
$ws.Range("E3").Value = 0.54000000000000004
$ws.Range("F3").Value = 47025.982386156596
$ws.Range("G3").Value = 0.99988063704640995
$ws.Range("H3").Value = 4370375.7890928704
$ws.Range("I3").Value = 71061.397815635399
$ws.Range("J3").Value = 92.935342704917304
# Row 4
$ws.Range("E4").Value = 0.55000000000000004
$ws.Range("F4").Value = 95602.023416329699
$ws.Range("G4").Value = 0.99952720551624097
$ws.Range("H4").Value = 4465916.1088990299
$ws.Range("I4").Value = 121177.81428218599
$ws.Range("J4").Value = 46.713614935227497
# Row 5
$ws.Range("E5").Value = 0.56000000000000005
$ws.Range("F5").Value = 147283.74390408001
$ws.Range("G5").Value = 0.99892451458078702
$ws.Range("H5").Value = 4564660.1610494098
$ws.Range("I5").Value = 174475.61194891299
$ws.Range("J5").Value = 30.992287675836
# Row 6
$ws.Range("E6").Value = 0.56999999999999995
$ws.Range("F6").Value = 202223.46142300399
$ws.Range("G6").Value = 0.998056741297348
$ws.Range("H6").Value = 4666757.4761447702
$ws.Range("I6").Value = 231109.83858068401
$ws.Range("J6").Value = 23.077230719451599
# Row 7
$ws.Range("E7").Value = 0.57999999999999996
$ws.Range("F7").Value = 260578.94007592701
$ws.Range("G7").Value = 0.99690729157643998
$ws.Range("H7").Value = 4772368.1341820098
$ws.Range("I7").Value = 291240.93753847602
$ws.Range("J7").Value = 18.314481334490999
# Row 8
$ws.Range("E8").Value = 0.58999999999999997
$ws.Range("F8").Value = 322513.27892964601
$ws.Range("G8").Value = 0.99545877334620503
$ws.Range("H8").Value = 4881663.9194660001
$ws.Range("I8").Value = 355034.60544504802
$ws.Range("J8").Value = 15.1363191483688
# Row 9
$ws.Range("E9").Value = 0.59999999999999998
$ws.Range("F9").Value = 388194.73720690003
$ws.Range("G9").Value = 0.99369297042475502
$ws.Range("H9").Value = 4994829.6523165302
$ws.Range("I9").Value = 422661.58128918801
$ws.Range("J9").Value = 12.8668144453859
# Row 10
$ws.Range("E10").Value = 0.60999999999999999
$ws.Range("F10").Value = 457796.48537665699
$ws.Range("G10").Value = 0.99159081744141497
$ws.Range("H10").Value = 5112064.7310608197
$ws.Range("I10").Value = 494297.35544947902
$ws.Range("J10").Value = 11.166675355436199
# Row 11
$ws.Range("E11").Value = 0.62
$ws.Range("F11").Value = 531496.26999659999
$ws.Range("G11").Value = 0.98913237620517602
$ws.Range("H11").Value = 5233584.9251558501
$ws.Range("I11").Value = 570121.78587727598
$ws.Range("J11").Value = 9.8468892833233603
# Row 12
$ws.Range("E12").Value = 0.63
$ws.Range("F12").Value = 609475.97886881197
$ws.Range("G12").Value = 0.98629681398210101
$ws.Range("H12").Value = 5359624.4694606801
$ws.Range("I12").Value = 650318.607339832
$ws.Range("J12").Value = 8.7938239656436394
# Row 13
$ws.Range("E13").Value = 0.64000000000000001
$ws.Range("F13").Value = 691921.09185484704
$ws.Range("G13").Value = 0.98306238421355896
$ws.Range("H13").Value = 5490438.5211862596
$ws.Range("I13").Value = 735074.81857884605
$ws.Range("J13").Value = 7.9350645410561604
# Row 14
$ws.Range("E14").Value = 0.65000000000000002
$ws.Range("F14").Value = 779020.00169389299
$ws.Range("G14").Value = 0.97940641028249897
$ws.Range("H14").Value = 5626306.05554647
$ws.Range("I14").Value = 824579.93137429794
$ws.Range("J14").Value = 7.2222870315430798
# Row 15
$ws.Range("E15").Value = 0.66000000000000003
$ws.Range("F15").Value = 870963.18859117304
$ws.Range("G15").Value = 0.97530527301458503
$ws.Range("H15").Value = 5767533.2944927104
$ws.Range("I15").Value = 919025.06529068504
$ws.Range("J15").Value = 6.6220172907904304
# Row 16
$ws.Range("E16").Value = 0.67000000000000004
$ws.Range("F16").Value = 967942.23252544401
$ws.Range("G16").Value = 0.97073440268169997
$ws.Range("H16").Value = 5914457.7862867098
$ws.Range("I16").Value = 1018601.87265792
$ws.Range("J16").Value = 6.1103416996853097
# Row 17
$ws.Range("E17").Value = 0.68000000000000005
$ws.Range("F17").Value = 1070148.6486785901
$ws.Range("G17").Value = 0.965668276352185
$ws.Range("H17").Value = 6067453.2835876103
$ws.Range("I17").Value = 1123501.28072289
$ws.Range("J17").Value = 5.6697294259817399
# Row 18
$ws.Range("E18").Value = 0.68999999999999995
$ws.Range("F18").Value = 1177772.5349117699
$ws.Range("G18").Value = 0.96008042149711603
$ws.Range("H18").Value = 6226935.6062606703
$ws.Range("I18").Value = 1233912.0429721801
$ws.Range("J18").Value = 5.2870443329934904
# Row 19
$ws.Range("E19").Value = 0.69999999999999996
$ws.Range("F19").Value = 1291001.0270248801
$ws.Range("G19").Value = 0.95394342680173505
$ws.Range("H19").Value = 6393369.7250369499
$ws.Range("I19").Value = 1350019.10081835
$ws.Range("J19").Value = 4.9522576599110097
# Row 20
$ws.Range("E20").Value = 0.70999999999999996
$ws.Range("F20").Value = 1410016.56953417
$ws.Range("G20").Value = 0.94722896112521304
$ws.Range("H20").Value = 6567278.3672513803
$ws.Range("I20").Value = 1472001.7726171301
$ws.Range("J20").Value = 4.6575894986971997
# Row 21
$ws.Range("E21").Value = 0.71999999999999997
$ws.Range("F21").Value = 1534995.02982774
$ws.Range("G21").Value = 0.93990780146752495
$ws.Range("H21").Value = 6749252.5313454596
$ws.Range("I21").Value = 1600031.81298244
$ws.Range("J21").Value = 4.3969214233240104
# Row 22
$ws.Range("E22").Value = 0.72999999999999998
$ws.Range("F22").Value = 1666103.7164108001
$ws.Range("G22").Value = 0.93194987058854495
$ws.Range("H22").Value = 6939964.4098033505
$ws.Range("I22").Value = 1734271.4274690901
$ws.Range("J22").Value = 4.1653855888118203
# Row 23
$ws.Range("E23").Value = 0.73999999999999999
$ws.Range("F23").Value = 1803499.4147817199
$ws.Range("G23").Value = 0.92332428450131299
$ws.Range("H23").Value = 7140183.3706623605
$ws.Range("I23").Value = 1874871.3950595399
$ws.Range("J23").Value = 3.9590716315960401
# Row 24
$ws.Range("E24").Value = 0.75
$ws.Range("F24").Value = 1947326.6388244601
$ws.Range("G24").Value = 0.91399940930333201
$ws.Range("H24").Value = 7350795.8496864801
$ws.Range("I24").Value = 2021969.5583659699
$ws.Range("J24").Value = 3.7748139953161202
# Row 25
$ws.Range("E25").Value = 0.76000000000000001
$ws.Range("F25").Value = 2097716.4300293601
$ws.Range("G25").Value = 0.90394292551699695
$ws.Range("H25").Value = 7572830.2784914998
$ws.Range("I25").Value = 2175690.1128270701
$ws.Range("J25").Value = 3.6100352602880301
# Row 26
$ws.Range("E26").Value = 0.77000000000000002
$ws.Range("F26").Value = 2254786.2516243602
$ws.Range("G26").Value = 0.89312189596528702
$ws.Range("H26").Value = 7807488.54656386
$ws.Range("I26").Value = 2336144.4008634798
$ws.Range("J26").Value = 3.4626291254611501
# Row 27
$ws.Range("E27").Value = 0.78000000000000003
$ws.Range("F27").Value = 2418641.8698159
$ws.Range("G27").Value = 0.88150282969763905
$ws.Range("H27").Value = 8056186.0077320598
$ws.Range("I27").Value = 2503434.3602353102
$ws.Range("J27").Value = 3.3308718038298402
# Row 28
$ws.Range("E28").Value = 0.79000000000000004
$ws.Range("F28").Value = 2589382.6737448899
$ws.Range("G28").Value = 0.869051728765038
$ws.Range("H28").Value = 8320602.75285096
$ws.Range("I28").Value = 2677660.49863191
$ws.Range("J28").Value = 3.2133538380471598
# Row 29
$ws.Range("E29").Value = 0.80000000000000004
$ws.Range("F29").Value = 2767112.8027830399
$ws.Range("G29").Value = 0.85573409533263201
$ws.Range("H29").Value = 8602749.86556191
$ws.Range("I29").Value = 2858937.4599935701
$ws.Range("J29").Value = 3.10892633538815
# Row 30
$ws.Range("E30").Value = 0.81000000000000005
$ws.Range("F30").Value = 2951961.9737730501
$ws.Range("G30").Value = 0.84151486137155995
$ws.Range("H30").Value = 8905055.7816444505
$ws.Range("I30").Value = 3047422.2483708402
$ws.Range("J30").Value = 3.0166566713128899
# Row 31
$ws.Range("E31").Value = 0.81999999999999995
$ws.Range("F31").Value = 3144122.4745338499
$ws.Range("G31").Value = 0.82635817797637301
$ws.Range("H31").Value = 9230479.8668263201
$ws.Range("I31").Value = 3243363.5867948001
$ws.Range("J31").Value = 2.93578890186039
# Row 32
$ws.Range("E32").Value = 0.82999999999999996
$ws.Range("F32").Value = 3343913.2174625299
$ws.Range("G32").Value = 0.81022695919078902
$ws.Range("H32").Value = 9582663.1774164103
$ws.Range("I32").Value = 3447186.8237684998
$ws.Range("J32").Value = 2.8657033105326999
# Row 33
$ws.Range("E33").Value = 0.83999999999999997
$ws.Range("F33").Value = 3551889.52458583
$ws.Range("G33").Value = 0.79308200354184599
$ws.Range("H33").Value = 9966130.4495193996
$ws.Range("I33").Value = 3659639.3584988099
$ws.Range("J33").Value = 2.8058672378560301
# Row 34
$ws.Range("E34").Value = 0.84999999999999998
$ws.Range("F34").Value = 3769031.3024376398
$ws.Range("G34").Value = 0.77488039228597805
$ws.Range("H34").Value = 10386563.1889208
$ws.Range("I34").Value = 3882040.8157841102
$ws.Range("J34").Value = 2.7557646396312001
# Row 35
$ws.Range("E35").Value = 0.85999999999999999
$ws.Range("F35").Value = 3997068.0880951099
$ws.Range("G35").Value = 0.75557264338328201
$ws.Range("H35").Value = 10851171.9298591
$ws.Range("I35").Value = 4116718.3275312702
$ws.Range("J35").Value = 2.71478285851029
# Row 36
$ws.Range("E36").Value = 0.87
$ws.Range("F36").Value = 4239048.5347681502
$ws.Range("G36").Value = 0.73509770055220203
$ws.Range("H36").Value = 11369206.822345899
$ws.Range("I36").Value = 4367777.1754512098
$ws.Range("J36").Value = 2.6820185541866399
# Row 37
$ws.Range("E37").Value = 0.88
$ws.Range("F37").Value = 4500358.2862722604
$ws.Range("G37").Value = 0.71337408944959602
$ws.Range("H37").Value = 11952659.292920399
$ws.Range("I37").Value = 4642497.0141726499
$ws.Range("J37").Value = 2.65593504618964
# Row 38
$ws.Range("E38").Value = 0.89000000000000001
$ws.Range("F38").Value = 4790586.3211697396
$ws.Range("G38").Value = 0.69028412958448904
$ws.Range("H38").Value = 12617219.5000793
$ws.Range("I38").Value = 4953934.9635768197
$ws.Range("J38").Value = 2.6337526670427498
# Row 39
$ws.Range("E39").Value = 0.90000000000000002
$ws.Range("F39").Value = 5127054.9305676799
$ws.Range("G39").Value = 0.66564519947120404
$ws.Range("H39").Value = 13383547.9537296
$ws.Range("I39").Value = 5325947.6101499796
$ws.Range("J39").Value = 2.6103773286953502
# Row 40
$ws.Range("E40").Value = 0.91000000000000003
$ws.Range("F40").Value = 5541744.9742929405
$ws.Range("G40").Value = 0.63915602431684204
$ws.Range("H40").Value = 14278836.880794
$ws.Range("I40").Value = 5803271.9347169204
$ws.Range("J40").Value = 2.5765958099895698
# Row 41
$ws.Range("E41").Value = 0.92000000000000004
$ws.Range("F41").Value = 6095463.1201239098
$ws.Range("G41").Value = 0.610292825126871
$ws.Range("H41").Value = 15338282.499053599
$ws.Range("I41").Value = 6465964.95236044
$ws.Range("J41").Value = 2.51634407374477
# Row 42
$ws.Range("E42").Value = 0.93000000000000005
$ws.Range("F42").Value = 6879816.5208035996
$ws.Range("G42").Value = 0.57813365070340705
$ws.Range("H42").Value = 16605756.4342955
$ws.Range("I42").Value = 7382301.4311472001
$ws.Range("J42").Value = 2.4136917582150601
# Row 43
$ws.Range("E43").Value = 0.93999999999999995
$ws.Range("F43").Value = 7929366.7573634004
$ws.Range("G43").Value = 0.54141697792215704
$ws.Range("H43").Value = 18142998.934820201
$ws.Range("I43").Value = 8563741.8952124808
$ws.Range("J43").Value = 2.2880766510102699
# Row 44
$ws.Range("E44").Value = 0.94999999999999996
$ws.Range("F44").Value = 9240221.2822580598
$ws.Range("G44").Value = 0.49882428350535202
$ws.Range("H44").Value = 20058850.342741702
$ws.Range("I44").Value = 9999608.8111673705
$ws.Range("J44").Value = 2.1708192617915198
# Row 45
$ws.Range("E45").Value = 0.95999999999999996
$ws.Range("F45").Value = 10803679.7878503
$ws.Range("G45").Value = 0.44909015003973002
$ws.Range("H45").Value = 22573660.725635301
$ws.Range("I45").Value = 11749420.645271201
$ws.Range("J45").Value = 2.08944185397104
# Row 46
$ws.Range("E46").Value = 0.96999999999999997
$ws.Range("F46").Value = 12835498.4069483
$ws.Range("G46").Value = 0.39065313859753398
$ws.Range("H46").Value = 26181740.752423301
$ws.Range("I46").Value = 14659637.8005598
$ws.Range("J46").Value = 2.0397915158674498
# Row 47
$ws.Range("E47").Value = 0.97999999999999998
$ws.Range("F47").Value = 17126518.847672898
$ws.Range("G47").Value = 0.31774184809619899
$ws.Range("H47").Value = 31942792.228355099
$ws.Range("I47").Value = 20991910.152890801
$ws.Range("J47").Value = 1.86510711910934
# Row 48
$ws.Range("E48").Value = 0.98999999999999999
$ws.Range("F48").Value = 25696294.365929101
$ws.Range("G48").Value = 0.21333631774422099
$ws.Range("H48").Value = 42893674.303819299
$ws.Range("I48").Value = 26485835.932297699
$ws.Range("J48").Value = 1.66925525108758
# Row 49
$ws.Range("E49").Value = 0.99099999999999999
$ws.Range("F49").Value = 27346578.898632601
$ws.Range("G49").Value = 0.20016330143847899
$ws.Range("H49").Value = 44716767.456210598
$ws.Range("I49").Value = 28339741.497626301
$ws.Range("J49").Value = 1.6351868956612501
# Row 50
$ws.Range("E50").Value = 0.99199999999999999
$ws.Range("F50").Value = 29382168.271236598
$ws.Range("G50").Value = 0.18606822519459501
$ws.Range("H50").Value = 46763895.701033697
$ws.Range("I50").Value = 30427916.459897701
$ws.Range("J50").Value = 1.5915740209960201
# Row 51
$ws.Range("E51").Value = 0.99299999999999999
$ws.Range("F51").Value = 31394489.003305402
$ws.Range("G51").Value = 0.170934572600032
$ws.Range("H51").Value = 49097607.021195903
$ws.Range("I51").Value = 32017479.261327799
$ws.Range("J51").Value = 1.56389253591694
# Row 52
$ws.Range("E52").Value = 0.99399999999999999
$ws.Range("F52").Value = 32554485.009273201
$ws.Range("G52").Value = 0.15501033379343601
$ws.Range("H52").Value = 51944294.981173903
$ws.Range("I52").Value = 33511489.768729001
$ws.Range("J52").Value = 1.59561101846266
# Row 53
$ws.Range("E53").Value = 0.995
$ws.Range("F53").Value = 34958545.812342003
$ws.Range("G53").Value = 0.13834303273933399
$ws.Range("H53").Value = 55630856.023662902
$ws.Range("I53").Value = 37485250.554207698
$ws.Range("J53").Value = 1.5913378182917099
# Row 54
$ws.Range("E54").Value = 0.996
$ws.Range("F54").Value = 40154661.817384697
$ws.Range("G54").Value = 0.11969933887829701
$ws.Range("H54").Value = 60167257.391026802
$ws.Range("I54").Value = 41437838.0783398
$ws.Range("J54").Value = 1.4983878500746799
# Row 55
$ws.Range("E55").Value = 0.997
$ws.Range("F55").Value = 42537687.094589502
$ws.Range("G55").Value = 0.099089782960334702
$ws.Range("H55").Value = 66410397.161922403
$ws.Range("I55").Value = 44670844.464345902
$ws.Range("J55").Value = 1.5612131664386899
# Row 56
$ws.Range("E56").Value = 0.998
$ws.Range("F56").Value = 50936663.654178597
$ws.Range("G56").Value = 0.076872256429389205
$ws.Range("H56").Value = 77280173.510710701
$ws.Range("I56").Value = 55649356.559737101
$ws.Range("J56").Value = 1.51718169127418
# Row 57
$ws.Range("E57").Value = 0.999
$ws.Range("F57").Value = 68480672.128838897
$ws.Range("G57").Value = 0.049194448440269503
$ws.Range("H57").Value = 98910990.461684301
$ws.Range("I57").Value = 69660810.049053907
$ws.Range("J57").Value = 1.4443636048956101
# Row 58
$ws.Range("E58").Value = 0.99909999999999999
$ws.Range("F58").Value = 70894227.709698707
$ws.Range("G58").Value = 0.045729792882510403
$ws.Range("H58").Value = 102161010.507532
$ws.Range("I58").Value = 72248341.070521802
$ws.Range("J58").Value = 1.44103425353423
# Row 59
$ws.Range("E59").Value = 0.99919999999999998
$ws.Range("F59").Value = 73670975.592775404
$ws.Range("G59").Value = 0.0421364436772316
$ws.Range("H59").Value = 105900094.18715701
$ws.Range("I59").Value = 75251258.983088598
$ws.Range("J59").Value = 1.43747375862826
# Row 60
$ws.Range("E60").Value = 0.99929999999999997
$ws.Range("F60").Value = 76922501.350863501
$ws.Range("G60").Value = 0.038393741109166102
$ws.Range("H60").Value = 110278499.216309
$ws.Range("I60").Value = 78807231.850801095
$ws.Range("J60").Value = 1.43363121686983
# Row 61
$ws.Range("E61").Value = 0.99939999999999996
$ws.Range("F61").Value = 80817741.580934793
$ws.Range("G61").Value = 0.034474178393326102
$ws.Range("H61").Value = 115523710.443893
$ws.Range("I61").Value = 83131427.051656306
$ws.Range("J61").Value = 1.4294350248354999
# Row 62
$ws.Range("E62").Value = 0.99950000000000006
$ws.Range("F62").Value = 85628824.297670901
$ws.Range("G62").Value = 0.030339547162995299
$ws.Range("H62").Value = 122002167.122339
$ws.Range("I62").Value = 88585879.712761194
$ws.Range("J62").Value = 1.4247791923222399
# Row 63
$ws.Range("E63").Value = 0.99960000000000004
$ws.Range("F63").Value = 91832791.183354899
$ws.Range("G63").Value = 0.025933632843743799
$ws.Range("H63").Value = 130356238.974731
$ws.Range("I63").Value = 95846929.2635528
$ws.Range("J63").Value = 1.4194955559442799

# New rows 64-75: full data (A/B/C blank, D blank, E-J values)
# Row 64
$ws.Range("E64").Value = 0.99970000000000003
$ws.Range("F64").Value = 100375316.651934
$ws.Range("G64").Value = 0.0211665823897852
$ws.Range("H64").Value = 141859342.211804
$ws.Range("I64").Value = 106400715.765119
$ws.Range("J64").Value = 1.4132891127379601
# Row 65
$ws.Range("E65").Value = 0.99980000000000002
$ws.Range("F65").Value = 113541599.40446199
$ws.Range("G65").Value = 0.015874627975738499
$ws.Range("H65").Value = 159588655.43514499
$ws.Range("I65").Value = 124647164.062047
$ws.Range("J65").Value = 1.40555229336388
# Row 66
$ws.Range("E66").Value = 0.99990000000000001
$ws.Range("F66").Value = 139490125.97274101
$ws.Range("G66").Value = 0.0096751667661666308
$ws.Range("H66").Value = 194530146.80824301
$ws.Range("I66").Value = 141624668.73390499
$ws.Range("J66").Value = 1.39458004967504
# Row 67
$ws.Range("E67").Value = 0.99990999999999997
$ws.Range("F67").Value = 143855579.79430601
$ws.Range("G67").Value = 0.0089707811937091899
$ws.Range("H67").Value = 200408533.260919
$ws.Range("I67").Value = 146304795.960154
$ws.Range("J67").Value = 1.3931231138025799
# Row 68
$ws.Range("E68").Value = 0.99992000000000003
$ws.Range("F68").Value = 148877947.93607101
$ws.Range("G68").Value = 0.0082431185034663196
$ws.Range("H68").Value = 207171500.42347601
$ws.Range("I68").Value = 151736243.06720099
$ws.Range("J68").Value = 1.39155263284819
# Row 69
$ws.Range("E69").Value = 0.99992999999999999
$ws.Range("F69").Value = 154759057.80537599
$ws.Range("G69").Value = 0.00748844192492881
$ws.Range("H69").Value = 215090822.90299499
$ws.Range("I69").Value = 158168013.467659
$ws.Range("J69").Value = 1.3898431920766301
# Row 70
$ws.Range("E70").Value = 0.99994000000000005
$ws.Range("F70").Value = 161804468.93372199
$ws.Range("G70").Value = 0.0067017762422701401
$ws.Range("H70").Value = 224577957.80883399
$ws.Range("I70").Value = 165989285.37641799
$ws.Range("J70").Value = 1.3879589314731799
# Row 71
$ws.Range("E71").Value = 0.99995000000000001
$ws.Range("F71").Value = 170506385.55765301
$ws.Range("G71").Value = 0.0058762106197820501
$ws.Range("H71").Value = 236295692.29539299
$ws.Range("I71").Value = 175854880.104359
$ws.Range("J71").Value = 1.3858465858775399
# Row 72
$ws.Range("E72").Value = 0.99995999999999996
$ws.Range("F72").Value = 181727644.005041
$ws.Range("G72").Value = 0.0050015773973420101
$ws.Range("H72").Value = 251405895.34307
$ws.Range("I72").Value = 188988108.42583501
$ws.Range("J72").Value = 1.3834213100572399
# Row 73
$ws.Range("E73").Value = 0.99997000000000003
$ws.Range("F73").Value = 197178707.46852201
$ws.Range("G73").Value = 0.0040616246466577601
$ws.Range("H73").Value = 272211824.315328
$ws.Range("I73").Value = 208076986.25434399
$ws.Range("J73").Value = 1.3805335667837499
# Row 74
$ws.Range("E74").Value = 0.99997999999999998
$ws.Range("F74").Value = 220992867.71166199
$ws.Range("G74").Value = 0.0030267313022291501
$ws.Range("H74").Value = 304279243.34599
$ws.Range("I74").Value = 241079759.45734799
$ws.Range("J74").Value = 1.3768735909748699
# Row 75
$ws.Range("E75").Value = 0.99999000000000005
$ws.Range("F75").Value = 267926568.086862
$ws.Range("G75").Value = 0.0018276951040063201
$ws.Range("H75").Value = 367478727.233926
$ws.Range("I75").Value = 367478727.23393703
$ws.Range("J75").Value = 1.3715650891134801
